$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in columns D/E hold numeric-looking or percent-looking text that must
# remain literal text (matching the source data exactly, including trailing
# zeros / thousand-dot groupings / padding spaces), so force Text format first.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.919.54'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.73%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.829.99'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.96%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9990'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.33'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.20%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6874'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -2.69%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9995'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07642'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -2.87%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3043'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -2.76%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.52'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -4.13%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07791'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -2.61%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.826.86'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -4.78%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.063'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.73%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '90.41'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -3.21%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6748'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -3.59%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.418'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.92%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008264'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.52%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '28.906.07'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -2.30%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '242.65'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -3.84%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.075.14'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -2.90%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.65'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -3.41%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9996'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.15%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.407'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -3.04%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.9993'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.10%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1471'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -5.70%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.24%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.763'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -2.93%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.17'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -2.98%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.534'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +2.30%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.210'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -2.66%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.117'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -3.69%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.82%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05106'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -3.72%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7482'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.50%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.823'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -3.54%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.143'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -2.34%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.672'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -1.32%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01842'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -2.02%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.213.86'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -4.64%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.682'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -2.22%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9168'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +2.24%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '108.19'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.83%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9989'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.06%  '
$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.496'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -9.59%  '
$ws.Range("B45").Value = 'Mantle'
$ws.Range("C45").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5165'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.25%  '
$ws.Range("B46").Value = 'RocketPoolETH'
$ws.Range("C46").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.974.51'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -2.81%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.501'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.82%  '
$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.00000000121'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -5.08%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '63.06'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -11.89%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.732'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -3.49%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4181'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -2.92%  '
